# Applies the "Handle case where a reasigned cajero has a perfect fit" edit.
#
# The change removes the "(T)" prefix from a handful of 1st-turno cajero
# names (those cells no longer represent a temporary/transferred worker)
# and re-shuffles several 2nd-turno schedule rows across sheets 1, 2, 4,
# 5, 6 and 7 so that a reassigned cajero who has a perfectly matching
# slot is placed correctly (sometimes leaving the row below empty).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5").Value  = "HEREDIA CAHUAYA, SUSAN NAYELLI"
$ws.Range("B6").Value  = "BRENIS LARTIGA, SEBASTIAN"

$ws.Range("F8").Value  = "LEON TICONA, MARIA FERNANDA"
$ws.Range("G8").Value  = "19:00-22:45"

$ws.Range("B15").Value = "ALVITE CORNEJO, ANGIE LUCERO"
$ws.Range("C15").Value = "10:15-14:00"
$ws.Range("F15").Value = "BARRIENTOS JERI, MILAGROS NICOL"
$ws.Range("G15").Value = "19:00-22:45"

$ws.Range("B16").Value = "SICHA JORGE, JOSE ANGELO"
$ws.Range("C16").Value = "15:15-19:00"
$ws.Range("F16").ClearContents()
$ws.Range("G16").ClearContents()

$ws.Range("B53").Value = "HUAMANI QUICANO, EMELYN HEIDY"
$ws.Range("B54").Value = "CARHUANCHO RAYMUNDO, JAMES JESUS"

# ---------------------------------------------------------------------
# Sheet 2
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Range("F5").Value  = "BRENIS LARTIGA, SEBASTIAN"

$ws.Range("B21").Value = "PIÑA SHUPINGAHUA, SANDRA"
$ws.Range("C21").Value = "18:15-22:00"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# ---------------------------------------------------------------------
# Sheet 4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)

$ws.Range("B5").Value  = "HUAMAN HUAMANI, ALEXIS JAVIER"

$ws.Range("B15").Value = "ALVITE CORNEJO, ANGIE LUCERO"
$ws.Range("C15").Value = "10:15-14:00"
$ws.Range("F15").ClearContents()
$ws.Range("G15").ClearContents()

$ws.Range("B16").Value = "HUAMANI LOPEZ, DREYDI BELINDA"
$ws.Range("C16").Value = "18:15-22:00"

# ---------------------------------------------------------------------
# Sheet 5
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)

$ws.Range("B5").Value  = "ERIQUE CALLE, MARIA ANTONIETA"
$ws.Range("B6").Value  = "RAMOS TINOCO, JORDAN JAVIER"

$ws.Range("B13").Value = "RUIZ SANTOS, CIELO CRISTHINA"
$ws.Range("C13").Value = "10:15-14:00"
$ws.Range("F13").Value = "HUAMANI LOPEZ, DREYDI BELINDA"
$ws.Range("G13").Value = "19:00-22:45"

$ws.Range("B14").Value = "VEGA CARDENAS, ANGELICA LOURDES"
$ws.Range("C14").Value = "15:00-18:45"
$ws.Range("F14").ClearContents()
$ws.Range("G14").ClearContents()

# ---------------------------------------------------------------------
# Sheet 6
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

$ws.Range("B5").Value  = "NINA OSCCO, FRANCIS GABRIELA"

$ws.Range("B17").Value = "RAMOS TINOCO, JORDAN JAVIER"
$ws.Range("C17").Value = "10:15-19:15"
$ws.Range("F17").ClearContents()
$ws.Range("G17").ClearContents()

$ws.Range("B18").Value = "VILCAPOMA CHILIN, JULISSA JAZMIN"
$ws.Range("C18").Value = "19:15-22:00"

# ---------------------------------------------------------------------
# Sheet 7
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(7)

$ws.Range("B5").Value  = "BRANCACHO RAMIREZ, BRINDY"

$ws.Range("B17").Value = "CASAPAICO RIVERA, ENZO MANUEL"
$ws.Range("C17").Value = "10:15-14:00"
$ws.Range("F17").ClearContents()
$ws.Range("G17").ClearContents()

$ws.Range("B18").Value = "RAMOS TINOCO, JORDAN JAVIER"
$ws.Range("C18").Value = "14:15-22:45"
